$d = $word.ActiveDocument

$replacements = @(
    @("2025-07-20 Sunday", "2025-07-21 Monday"),
    @("77-61=16", "20+1=21"),
    @("46-9=37", "53-45=8"),
    @("60+39=99", "63+36=99"),
    @("59+15=74", "96-89=7"),
    @("77-28=49", "37+2=39"),
    @("50-9=41", "82-8=74"),
    @("48-33=15", "52-18=34"),
    @("93-53=40", "37+52=89"),
    @("52+18=70", "11+32=43"),
    @("38-25=13", "86-54=32"),
    @("19+9=28", "3+84=87"),
    @("84-80=4", "33+42=75"),
    @("83-55=28", "20+50=70"),
    @("43+49=92", "95-71=24"),
    @("20+53=73", "48+33=81"),
    @("73+5=78", "94-1=93"),
    @("80-21=59", "15+53=68"),
    @("85-38=47", "83+13=96"),
    @("12+70=82", "30-8=22"),
    @("91-67=24", "18+42=60"),
    @("46-2=44", "0+10=10"),
    @("14+49=63", "84-31=53"),
    @("25+71=96", "76+2=78"),
    @("64-23=41", "1+62=63"),
    @("62+36=98", "85+13=98"),
    @("6-1=5", "54-15=39"),
    @("98-5=93", "94-13=81"),
    @("35+26=61", "60+27=87"),
    @("15+25=40", "89-17=72"),
    @("98-63=35", "59-53=6"),
    @("1-0=1", "65-55=10"),
    @("21+51=72", "9+75=84"),
    @("34+61=95", "44-11=33"),
    @("41-26=15", "78+19=97"),
    @("46-37=9", "36-18=18"),
    @("56-44=12", "26+5=31"),
    @("28-20=8", "15+74=89"),
    @("74+11=85", "62-62=0"),
    @("73+21=94", "64+3=67"),
    @("83-24=59", "28-1=27"),
    @("31-16=15", "79-15=64"),
    @("14-0=14", "32+52=84"),
    @("64-34=30", "35+57=92"),
    @("15+70=85", "60+27=87"),
    @("92-80=12", "48+49=97"),
    @("60-56=4", "81+14=95"),
    @("12+51=63", "43-4=39"),
    @("77+21=98", "73-60=13"),
    @("21-3=18", "42-35=7"),
    @("49+20=69", "5+50=55"),
    @("88-69=19", "48+7=55"),
    @("79+3=82", "78-73=5"),
    @("20+65=85", "21+21=42"),
    @("40-25=15", "77-54=23"),
    @("46-7=39", "22-2=20"),
    @("50-3=47", "5+43=48"),
    @("61-12=49", "12+83=95"),
    @("14-6=8", "80-71=9"),
    @("54+16=70", "0+56=56"),
    @("92+3=95", "23-13=10"),
    @("14-12=2", "96-73=23"),
    @("5+63=68", "28+6=34"),
    @("56+11=67", "8+17=25"),
    @("88-21=67", "34+38=72"),
    @("83-61=22", "83-2=81"),
    @("99-13=86", "29+69=98"),
    @("29-23=6", "50+40=90"),
    @("7+80=87", "28-16=12"),
    @("42+0=42", "35+45=80"),
    @("57+30=87", "40+57=97"),
    @("2+30=32", "91-3=88"),
    @("91-81=10", "28-15=13"),
    @("56-41=15", "87-70=17"),
    @("88-3=85", "92-8=84"),
    @("76+1=77", "1+38=39"),
    @("19+65=84", "95-14=81"),
    @("49-10=39", "14+78=92"),
    @("19+28=47", "44-37=7"),
    @("19+57=76", "90-87=3"),
    @("97-10=87", "52-12=40"),
    @("3+17=20", "1+87=88"),
    @("26+45=71", "94-94=0"),
    @("56+40=96", "3+57=60"),
    @("44-24=20", "98-18=80"),
    @("88-7=81", "95-31=64"),
    @("96-58=38", "65-43=22"),
    @("36+43=79", "91-68=23"),
    @("3+82=85", "86-47=39"),
    @("87-31=56", "51-31=20"),
    @("34+63=97", "71-30=41"),
    @("8+91=99", "60-37=23"),
    @("34+46=80", "0+94=94"),
    @("75-61=14", "74-2=72"),
    @("43+9=52", "37+13=50"),
    @("32+50=82", "14+43=57"),
    @("65+19=84", "61+27=88"),
    @("65-41=24", "49+45=94"),
    @("30+51=81", "96-44=52"),
    @("4+26=30", "55+20=75"),
    @("65-56=9", "25+31=56"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
